$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A3 was the text "asw" (a stray test value); replace it with the
# numeric value 1 to match the rest of column A.
$ws.Range("A3").Value = 1

# Move the active selection from D3 to D2.
$ws.Range("D2").Select()
